$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 data values (B2:E2)
$ws.Range("B2").Value = 4.2471599082983564
$ws.Range("C2").Value = 7.583308857097931
$ws.Range("D2").Value = 10.755305864937521
$ws.Range("E2").Value = 9.6949890951225068

# Update row 3 data values (B3:E3)
$ws.Range("B3").Value = 5.2817447298443208
$ws.Range("C3").Value = 6.1480008132484265
$ws.Range("D3").Value = 5.420649088666738
$ws.Range("E3").Value = 10.878692251859331

# Update the selection to match the new narrower highlighted range
$ws.Range("B1:E3").Select()
